$d = $word.ActiveDocument

# The timesheet table's last row (date "1/15") has two empty cells that
# need to be filled in with the final log entries for the week.
$table = $d.Tables.Item(1)
$lastRow = $table.Rows.Count

$timeCell = $table.Cell($lastRow, 2)
$timeCell.Range.Text = "9:05 – 9:50 AM`r2:50 – 3:00 PM"

$notesCell = $table.Cell($lastRow, 3)
$notesCell.Range.Text = "Checked answers for labs.`rCompleted Heaps labs."
